# Weekly data update: a new "Puerro" (leek) price observation is inserted
# as row 168 (date 2022-06-24 / serial 44736), pushing the existing rows
# 168-213 down to 169-214 (the sheet's dimension grows from R213 to R214).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 168; Excel shifts rows 168:213 down to 169:214
# and the sheet dimension updates automatically to A1:R214.
$ws.Rows("168:168").Insert()

# Populate the newly inserted row 168 with the new observation.
$ws.Cells.Item(168, 1).Value2 = 10
$ws.Cells.Item(168, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(168, 3).Value2 = "La Araucanía"
$ws.Cells.Item(168, 4).Value2 = 44736
$ws.Cells.Item(168, 5).Value2 = 9
$ws.Cells.Item(168, 6).Value2 = 100112005
$ws.Cells.Item(168, 7).Value2 = "Puerro"
$ws.Cells.Item(168, 8).Value2 = "Azul de Maquehue"
$ws.Cells.Item(168, 9).Value2 = "Primera"
$ws.Cells.Item(168, 10).Value2 = 30
$ws.Cells.Item(168, 11).Value2 = 15000
$ws.Cells.Item(168, 12).Value2 = 15000
$ws.Cells.Item(168, 13).Value2 = 15000
$ws.Cells.Item(168, 14).Value2 = "`$/docena de paquetes"
$ws.Cells.Item(168, 15).Value2 = "Región de La Araucanía"
$ws.Cells.Item(168, 16).Value2 = 1250
$ws.Cells.Item(168, 17).Value2 = 12
$ws.Cells.Item(168, 18).Value2 = "Hortaliza"
